$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the confidence interval value for L2 (was computed incorrectly)
$ws.Range("L2").Value = 1.29

# Remove the now-unneeded "gs_se" column (M) entirely
$ws.Columns.Item(13).Delete()
